$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (date) for rows 2 through 89 from 2023-09-23 (45192) to 2023-10-03 (45202)
$ws.Range("C2:C89").Value = 45202
